# move all things to sb_sbsj KJS_Planner branch
# - rename the "id" PK fields to their real column names (planId / planSpId)
# - add a new memo row (JPA를 잊어라 / Mybatis 사용)
# - rezoom + reselect the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# planner_tbl's id column -> planId
$ws.Range("A7").Value = "planId"

# planner_spot_tbl's id column -> planSpId
$ws.Range("A20").Value = "planSpId"

# new note row under the first table
$ws.Range("B15").Value = "JPA를 잊어라"
$ws.Range("C15").Value = "Mybatis 사용"

# view changes: zoom out to 70% and move the selection
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("D15").Select()
